$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Geanne.S83@gmail.com", 63014224),
    @("B_Mons32@yahoo.com", 32992576),
    @("SaadetB29@kpnmail.nl", 56928796),
    @("Ilian_B@gmail.com", 93898430),
    @("L.Lagendijk@live.com", 38370214),
    @("Sippie_O@hotmail.com", 44864502),
    @("S.Hameleers26@kpnmail.nl", 16331942)
)

$startRow = 31
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
